$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume updates ---
$ws.Range("D2").Value = "27.216.39"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.688.01"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +13.43%  "
$ws.Range("E9").Value = "  +3.76%  "
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0891"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.925.30"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "1.699.47"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "27.205.32"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").Value = "1.582.41"
$ws.Range("E32").Value = "  +6.70%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("E36").Value = "  +5.96%  "
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D45").Value = "1.833.79"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("E48").Value = "  +5.77%  "
$ws.Range("E49").Value = "  +3.18%  "

# --- Row swaps (coin order changed) with updated values ---
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.47%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
